$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 18 new rows above the existing "syngas, RWGS, Production" block (old row 217),
#    pushing everything from old row 217 down to row 235.
$ws.Rows("217:234").Insert()

# 2. New activity header block (rows 217-224)
$ws.Range("A217").Value = "Activity"
$ws.Range("A217").Font.Bold = $true
$ws.Range("A217").Font.Size = 12
$ws.Range("B217").Value = "liquefied petroleum gas production, synthetic, Fischer Tropsch process, hydrogen from wood gasification, energy allocation, with carbon capture and storage"
$ws.Range("B217").Font.Bold = $true

$ws.Range("A218").Value = "production amount"
$ws.Range("B218").Value = 1

$ws.Range("A219").Value = "reference product"
$ws.Range("B219").Value = "liquefied petroleum gas, synthetic"

$ws.Range("A220").Value = "type"
$ws.Range("B220").Value = "process"

$ws.Range("A221").Value = "unit"
$ws.Range("B221").Value = "kilogram"

$ws.Range("A222").Value = "location"
$ws.Range("B222").Value = "RER"

$ws.Range("A223").Value = "comment"
$ws.Range("B223").Value = "Adapted from A.E.M. van den Oever, D. Costa, M. Messagie, Prospective life cycle assessment of alternatively fueled heavy-duty trucks, Applied Energy, 2023, https://doi.org/10.1016/j.apenergy.2023.120834. Energy-based allocation between 3 other co-products (diesel, wax and C5-C10 olefins). Allocation key for this co-product: 11%. Post corrected to preserve carbon balance."

$ws.Range("A224").Value = "source"
$ws.Range("B224").Value = " A.E.M. van den Oever, D. Costa, M. Messagie, Prospective life cycle assessment of alternatively fueled heavy-duty trucks, Applied Energy, 2023, https://doi.org/10.1016/j.apenergy.2023.120834"

# 3. Exchanges header (row 225) + column titles (row 226)
$ws.Range("A225").Value = "Exchanges"
$ws.Range("A225").Font.Bold = $true
$ws.Range("A225").Font.Size = 12

$ws.Range("A226").Value = "name"
$ws.Range("B226").Value = "amount"
$ws.Range("C226").Value = "location"
$ws.Range("D226").Value = "unit"
$ws.Range("E226").Value = "categories"
$ws.Range("F226").Value = "type"
$ws.Range("G226").Value = "reference product"
$ws.Range("H226").Value = "comment"

# 4. Exchange rows (227-233)
# row 227 - reference product (self-reference via formula, like other blocks)
$ws.Range("A227").Formula = "=B217"
$ws.Range("B227").Value = 1
$ws.Range("C227").Value = "RER"
$ws.Range("D227").Value = "kilogram"
$ws.Range("F227").Value = "production"
$ws.Range("G227").Formula = "=B217"
$ws.Range("H227").Value = "From PEM electrolysis"

# row 228 - syngas input
$ws.Range("A228").Value = "syngas, RWGS, Production, for Fischer Tropsch process, hydrogen from wood gasification, with CCS"
$ws.Range("B228").Value = 2.28
$ws.Range("C228").Value = "RER"
$ws.Range("D228").Value = "kilogram"
$ws.Range("F228").Value = "technosphere"
$ws.Range("G228").Value = "syngas, RWGS, Production"

# row 229 - waste heat
$ws.Range("A229").Value = "Heat, waste"
$ws.Range("B229").Formula = "=3.64160231884058*0.11"
$ws.Range("D229").Value = "megajoule"
$ws.Range("E229").Value = "air"
$ws.Range("F229").Value = "biosphere"

# row 230 - water
$ws.Range("A230").Value = "Water"
$ws.Range("B230").Formula = "=0.00107549913043478*0.11"
$ws.Range("D230").Value = "cubic meter"
$ws.Range("E230").Value = "water"
$ws.Range("F230").Value = "biosphere"

# row 231 - gas-to-liquid plant construction
$ws.Range("A231").Value = "Gas-to-liquid plant construction"
$ws.Range("B231").Value = 0.0000000000067
$ws.Range("B231").NumberFormat = "0.00E+00"
$ws.Range("C231").Value = "GLO"
$ws.Range("D231").Value = "unit"
$ws.Range("F231").Value = "technosphere"
$ws.Range("G231").Value = "Gas-to-liquid plant"

# row 232 - CO2 uptake
$ws.Range("A232").Value = "Carbon dioxide, in air"
$ws.Range("B232").Value = -0.13
$ws.Range("D232").Value = "kilogram"
$ws.Range("E232").Value = "natural resource::in air"
$ws.Range("F232").Value = "biosphere"

# row 233 - electricity
$ws.Range("A233").Value = "market group for electricity, low voltage"
$ws.Range("B233").Formula = "=0.0872420618556701*0.11"
$ws.Range("B233").NumberFormat = "0.00E+00"
$ws.Range("C233").Value = "RER"
$ws.Range("D233").Value = "kilowatt hour"
$ws.Range("F233").Value = "technosphere"
$ws.Range("G233").Value = "electricity, low voltage"

# 5. Refresh the AutoFilter range to cover the new dimension
$ws.AutoFilterMode = $false
$ws.Range("A1:K596").AutoFilter()

# 6. Update the hidden _FilterDatabase defined name to match
foreach ($n in $wb.Names) {
    if ($n.Name -eq "FT fuel - Diesel!_FilterDatabase") {
        $n.RefersTo = "='FT fuel - Diesel'!`$A`$1:`$K`$596"
    }
}

# 7. Restore the view/selection state to match the edited workbook
$ws.Range("B163").Select()
